$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '22.320.96'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +8.60%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.586.33'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +7.71%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.55%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.9919'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '298.72'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +7.66%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3601'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3325'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +7.79%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '40.79'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.108'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.03%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06889'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.88%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9984'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.22%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.22'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +5.63%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.764'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.10%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.461'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.79%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9929'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.40%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001057'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.23%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.591.87'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +8.19%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06560'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +10.04%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '75.81'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +10.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.74'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +7.93%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.864'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +6.94%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.40'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '22.322.05'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +8.62%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.365'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.19%  '
$ws.Range('B26').Value = 'LEO'
$ws.Range('C26').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.400'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -7.16%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.487'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +17.29%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '148.53'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.12%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.01'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +10.87%  '
$ws.Range('B30').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C30').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.763.96'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +8.06%  '
$ws.Range('B31').Value = 'BitcoinCash'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '122.25'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +7.39%  '
$ws.Range('B32').Value = 'HuobiToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.922'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.01%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.805'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +17.47%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9131'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +13.10%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.08073'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.91%  '
$ws.Range('B36').Value = 'WEMIXTOKEN'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.625'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +10.49%  '
$ws.Range('B37').Value = 'Aptos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '11.61'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +11.68%  '
$ws.Range('E38').Value = '  +0.30%  '
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.033'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.56%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.327'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +12.55%  '
$ws.Range('B41').Value = 'Hedera'
$ws.Range('C41').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.05962'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.88%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.02167'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +5.72%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9910'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.06%  '
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1962'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.75%  '
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5729'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +8.71%  '
$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.747'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.54%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.57'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.88%  '
$ws.Range('B48').Value = 'Decentraland'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5556'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +6.87%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '123.16'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.39%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.920'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.65%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06733'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.39%  '
